$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force changed Price (D column) cells to Text format first so that
# numeric-looking strings (e.g. "609.28") are preserved exactly as text
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.256.77"
$ws.Range("D3").Value = "3.138.37"
$ws.Range("D5").Value = "609.28"
$ws.Range("D6").Value = "143.37"
$ws.Range("D8").Value = "3.135.07"
$ws.Range("D11").Value = "5.41"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("D14").Value = "35.51"
$ws.Range("D15").Value = "3.656.20"
$ws.Range("D17").Value = "64.247.10"
$ws.Range("D18").Value = "3.146.69"
$ws.Range("D19").Value = "6.87"
$ws.Range("D20").Value = "476.85"
$ws.Range("D21").Value = "14.71"
$ws.Range("D22").Value = "0.723"
$ws.Range("D23").Value = "7.82"
$ws.Range("D24").Value = "13.69"
$ws.Range("D25").Value = "85.39"
$ws.Range("D27").Value = "2.78"
$ws.Range("D28").Value = "8.55"
$ws.Range("D29").Value = "7.37"
$ws.Range("D31").Value = "2.08"
$ws.Range("D33").Value = "26.65"
$ws.Range("D35").Value = "1.10"
$ws.Range("D36").Value = "5.95"
$ws.Range("D37").Value = "52.40"
$ws.Range("D38").Value = "0.0₃0745"
$ws.Range("D39").Value = "455.24"
$ws.Range("D43").Value = "8.33"
$ws.Range("D44").Value = "2.874.58"
$ws.Range("D45").Value = "0.263"
$ws.Range("D46").Value = "2.26"
$ws.Range("D48").Value = "26.51"
$ws.Range("D51").Value = "120.96"

# Restore default cell style (keeps workbook styles clean) now that the
# text values are committed.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"

$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  +8.61%  "
$ws.Range("E30").Value = "  +3.13%  "
$ws.Range("E31").Value = "  -4.83%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("E38").Value = "  +4.61%  "
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  +6.16%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +2.45%  "
